{"js": "// Find the specific paragraph that contains the \"Defendant appeared in Court for\"\n// template text (it has a tab stop at pos 4680, centered) and set its alignment\n// to Justify (OOXML w:jc val=\"both\"), matching the author's \"clean format\" edit.\nconst searchResults = context.document.body.search(\"Defendant appeared in Court for\", { matchCase: false, matchWholeWord: false });\nsearchResults.load(\"paragraphs\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  const range = searchResults.items[i];\n  const paragraphs = range.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  for (let j = 0; j < paragraphs.items.length; j++) {\n    paragraphs.items[j].alignment = Word.Alignment.justified;\n  }\n}\nawait context.sync();\n", "ps1": "# Locate the specific paragraph that contains the \"Defendant appeared in Court for\"\n# template text (it has a tab stop at pos 4680) and set its alignment to\n# Justify (OOXML w:jc val=\"both\"), matching the author's \"clean format\" edit.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"Defendant appeared in Court for\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 1  # wdFindContinue\n$range.Find.MatchWildcards = $false\n\nwhile ($range.Find.Execute()) {\n    $range.Paragraphs.Alignment = 3  # wdAlignParagraphJustify\n    $range.Collapse(0)  # wdCollapseEnd\n}\n"}
